# "Generate Report for Handback"
#
# The handback round-trip for both target locales (zh-cn, de-de) has now
# completed. This script updates the localization-status workbook to
# reflect that:
#   - The Overview sheet's per-locale status changes from
#     "Ready for handoff" to "Handed back: in sync with en-US".
#   - Each locale sheet (zh-cn / de-de) gets its "Latest Target File",
#     "Latest Handback File" and "Latest Handback DateTime" columns
#     filled in for both source files, with a hyperlink added on the
#     "Latest Target File" cell (mirroring the existing hyperlink on the
#     "Source File Name" cell).
#   - A handful of columns are widened so the new, longer values are not
#     clipped.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
# everywhere it is used - the Overview sheet's per-locale columns as well
# as the "Status" column on each locale sheet.
# ---------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# Widen the now-longer zh-cn / de-de status columns on the Overview sheet.
$wsOverview.Range("E1:F1").ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# Locale sheets (zh-cn, de-de): fill in target/handback info for the two
# source files (rows 2 and 3).
# ---------------------------------------------------------------------

function Fill-HandbackRow {
    param(
        $ws,
        [int]$row,
        [string]$targetDisplay,
        [string]$targetUrl,
        [string]$handbackFile,
        [string]$handbackDateTime
    )

    # "Latest Target File" (column I) - hyperlinked, like column A.
    $ws.Hyperlinks.Add($ws.Range("I" + $row), $targetUrl, "", "", $targetDisplay)

    # "Latest Handback File" (column J)
    $ws.Range("J" + $row).Value = $handbackFile

    # "Latest Handback DateTime" (column K)
    $ws.Range("K" + $row).Value = $handbackDateTime
}

$url1aa = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/58e90c183d990db36743bb915e13c91c591a1d74/e2e/1aa29009-39e0-4b33-a645-3f348e20e891.md"
$url603 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/58e90c183d990db36743bb915e13c91c591a1d74/e2e/603718cb-1111-4a69-ba0a-989b0d347a7d.md"

# zh-cn
Fill-HandbackRow $wsZhCn 2 "1aa29009-39e0-4b33-a645-3f348e20e891.md" $url1aa `
    "1aa29009-39e0-4b33-a645-3f348e20e891.1d419a78037f0c5f01dfa176c821250c8473c753.zh-cn.xlf" `
    "2016-08-26 17:02:36"
Fill-HandbackRow $wsZhCn 3 "603718cb-1111-4a69-ba0a-989b0d347a7d.md" $url603 `
    "603718cb-1111-4a69-ba0a-989b0d347a7d.57328d7613f7bf05c785a2af73361c52d54d9c34.zh-cn.xlf" `
    "2016-08-26 17:02:36"

# de-de
Fill-HandbackRow $wsDeDe 2 "1aa29009-39e0-4b33-a645-3f348e20e891.md" $url1aa `
    "1aa29009-39e0-4b33-a645-3f348e20e891.1d419a78037f0c5f01dfa176c821250c8473c753.de-de.xlf" `
    "2016-08-26 17:02:43"
Fill-HandbackRow $wsDeDe 3 "603718cb-1111-4a69-ba0a-989b0d347a7d.md" $url603 `
    "603718cb-1111-4a69-ba0a-989b0d347a7d.57328d7613f7bf05c785a2af73361c52d54d9c34.de-de.xlf" `
    "2016-08-26 17:02:43"

# Widen the Status / Latest Target File / Latest Handback File columns on
# both locale sheets so the new values display fully.
$wsZhCn.Range("C1").ColumnWidth = 29.166666666666668
$wsZhCn.Range("I1:J1").ColumnWidth = 39.166666666666664

$wsDeDe.Range("C1").ColumnWidth = 29.166666666666668
$wsDeDe.Range("I1:J1").ColumnWidth = 39.166666666666664
